# TP FINAL ARREGLADO Y TERMINADO
#
# 1) Remove the empty "Hoja1" sheet (it had no data, sat between
#    "Dominios" and "Tablas").
# 2) Add the trigger / trigger-description pair ("Triggers de la tabla" /
#    "Descripcion2" columns, G/H) documenting the new "pasar a mayusculas"
#    triggers for idioma, ciudad, provincia, pais, autor, tema, editorial
#    and usuario/correo tables.
# 3) Leave "Tablas" as the selected/active sheet, with G61 selected.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# --- 1) drop the unused empty sheet -----------------------------------
$hoja1 = $wb.Worksheets.Item("Hoja1")
[void]$hoja1.Delete()

# --- 2) document the new "minAMayus*" triggers on the Tablas sheet ----
$ws = $wb.Worksheets.Item("Tablas")
$ws.Activate()

# idioma
$ws.Range("G17").Value = "minAMayus"
$ws.Range("H17").Value = "Pasa el nombre del idioma a mayusculas"

# ciudad
$ws.Range("G12").Value = "minAMayus"
$ws.Range("H12").Value = "Pasa el nombre de la ciudad a mayusculas"

# provincia
$ws.Range("G7").Value = "minAMayus"
$ws.Range("H7").Value = "Pasa el nombre de la provincia a mayusculas"

# pais
$ws.Range("G3").Value = "minAMayus"
$ws.Range("H3").Value = "Pasa el nombre del pais a mayusculas"

# autor
$ws.Range("G21").Value = "minAMayus"
$ws.Range("H21").Value = "Pasa el nombre del autor a mayusculas"

# tema
$ws.Range("G25").Value = "minAMayus"
$ws.Range("H25").Value = "Pasa el nombre del tema a mayusculas"

# editorial
$ws.Range("G29").Value = "minAMayus"
$ws.Range("H29").Value = "Pasa el nombre de la editorial a mayusculas"

# usuario / correo
$ws.Range("G60").Value = "minAMayusCorreo"
$ws.Range("H60").Value = "Pasa el correo a mayusculas"

# --- 3) leave the view the way the author left it ----------------------
[void]$ws.Range("G61").Select()
